$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("B12")
$r.Borders.Item(8).Color = 14540253
$r.Borders.Color = 15658734
Write-Host "done"
